# Insert a new data row at row 103 (this shifts existing rows 103:191 down to 104:192,
# and carries formatting along automatically).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(103).Insert()

# Fill in the new record for row 103.
$ws.Cells.Item(103, 1).Value = 10
$ws.Cells.Item(103, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(103, 3).Value = "La Araucanía"
$ws.Cells.Item(103, 4).Value = [datetime]"2022-11-09"
$ws.Cells.Item(103, 5).Value = 9
$ws.Cells.Item(103, 6).Value = "Fruta"
$ws.Cells.Item(103, 7).Value = 100104
$ws.Cells.Item(103, 8).Value = "Frutos de pepita"
$ws.Cells.Item(103, 9).Value = 100104001
$ws.Cells.Item(103, 10).Value = "Granada"
$ws.Cells.Item(103, 11).Value = "Wonderfull"
$ws.Cells.Item(103, 12).Value = "Primera"
$ws.Cells.Item(103, 13).Value = 80
$ws.Cells.Item(103, 14).Value = 16000
$ws.Cells.Item(103, 15).Value = 16000
$ws.Cells.Item(103, 16).Value = 16000
$ws.Cells.Item(103, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(103, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(103, 19).Value = 1067
$ws.Cells.Item(103, 20).Value = 15
